$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3-17, keeping only the header row and the first data row.
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row with the new values.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2.061751933828537
